# Update latest output (run 103)
$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item(1)
$wsDetailed = $wb.Worksheets.Item(2)

$dtFmt = "YYYY-MM-DD HH:MM:SS"
$dFmt = "YYYY-MM-DD"

# ---------------------------------------------------------------
# Sheet 1 "Schedule": replace the single data row with four rows
# Columns: A=Start Time, B=Stop Time, C=Duration(h), D=Volume(ML),
#          E=Cost($), F=Unit Cost($/ML)
# ---------------------------------------------------------------
$scheduleData = @"
2|46041|46041.22916666666|5.5|20.79|509.8020524999999|24.52150324675324
3|46041.27083333334|46041.66666666666|9.5|35.91|123.0863205|3.427633542188806
4|46041.95833333334|46042.125|4|15.12|456.29010375|30.17791691468254
5|46042.29166666666|46042.66666666666|9|34.02|-32.1935055|-0.9463111552028221
"@

$scheduleLines = $scheduleData -split "`n"
foreach ($line in $scheduleLines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $p = $line.Split("|")
    $r = [int]$p[0]

    $wsSchedule.Cells.Item($r, 1).Value = [double]$p[1]
    $wsSchedule.Cells.Item($r, 1).NumberFormat = $dtFmt

    $wsSchedule.Cells.Item($r, 2).Value = [double]$p[2]
    $wsSchedule.Cells.Item($r, 2).NumberFormat = $dtFmt

    $wsSchedule.Cells.Item($r, 3).Value = [double]$p[3]
    $wsSchedule.Cells.Item($r, 4).Value = [double]$p[4]
    $wsSchedule.Cells.Item($r, 5).Value = [double]$p[5]
    $wsSchedule.Cells.Item($r, 6).Value = [double]$p[6]
}

# ---------------------------------------------------------------
# Sheet 2 "Detailed": update existing rows 13,14 and 25-49
# Columns: A=DateTime, B=Price, C=Type, D=Date, E=Pump_Status
# Pipe-delimited: row|newPrice|newType|newPumpStatus (blank = unchanged)
# ---------------------------------------------------------------
$detailedChanges = @"
13|||OFF
14|||OFF
25|-4.82586||
26|-5.50985||
27|-5.71383|historical|
28|-5.42612|historical|
29|4.83168|historical|
30|13.72616|historical|
31|13.52945|historical|
32|22.76406||
33|12.69671||
34|12.09266||
35|0||
36|0.25443||
37|7.24587||
38|12.90677||
39|37.04015||
40|59.10107||
41|65||
42|65||
43|64.38026000000001||
44|63.12824||
45|64.31842||
46|59.27172||
47|59.34902||
48|58.32029||ON
49|62.01845||ON
"@

$detailedLines = $detailedChanges -split "`n"
foreach ($line in $detailedLines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $p = $line.Split("|")
    $r = [int]$p[0]
    $priceStr = $p[1]
    $typeStr = $p[2]
    $statusStr = $p[3]

    if ($priceStr -ne "") {
        $wsDetailed.Cells.Item($r, 2).Value = [double]$priceStr
    }
    if ($typeStr -ne "") {
        $wsDetailed.Cells.Item($r, 3).Value = $typeStr
    }
    if ($statusStr -ne "") {
        $wsDetailed.Cells.Item($r, 5).Value = $statusStr
    }
}

# ---------------------------------------------------------------
# Sheet 2 "Detailed": append new rows 50-97
# Pipe-delimited: row|DateTime|Price|Type|Date|Pump_Status
# ---------------------------------------------------------------
$newRows = @"
50|46042|62.41679|forecast|46042|ON
51|46042.02083333334|57.06003|forecast|46042|ON
52|46042.04166666666|57.06003|forecast|46042|ON
53|46042.0625|57.06003|forecast|46042|ON
54|46042.08333333334|57.06003|forecast|46042|ON
55|46042.10416666666|56.9942|forecast|46042|ON
56|46042.125|57.06003|forecast|46042|OFF
57|46042.14583333334|57.06003|forecast|46042|OFF
58|46042.16666666666|57.06003|forecast|46042|OFF
59|46042.1875|57.06003|forecast|46042|OFF
60|46042.20833333334|64.89|forecast|46042|OFF
61|46042.22916666666|67.38887|forecast|46042|OFF
62|46042.25|75.00427999999999|forecast|46042|OFF
63|46042.27083333334|70.00581|forecast|46042|OFF
64|46042.29166666666|36.06|forecast|46042|ON
65|46042.3125|8.717219999999999|forecast|46042|ON
66|46042.33333333334|4.48583|forecast|46042|ON
67|46042.35416666666|8.43693|forecast|46042|ON
68|46042.375|0.7|forecast|46042|ON
69|46042.39583333334|0.00025|forecast|46042|ON
70|46042.41666666666|-2.54265|forecast|46042|ON
71|46042.4375|-6.55839|forecast|46042|ON
72|46042.45833333334|-8.09287|forecast|46042|ON
73|46042.47916666666|-9.31202|forecast|46042|ON
74|46042.5|-8.360110000000001|forecast|46042|ON
75|46042.52083333334|-9.593819999999999|forecast|46042|ON
76|46042.54166666666|-7.97797|forecast|46042|ON
77|46042.5625|-8.651149999999999|forecast|46042|ON
78|46042.58333333334|-8.426119999999999|forecast|46042|ON
79|46042.60416666666|-7.76554|forecast|46042|ON
80|46042.625|-7.70186|forecast|46042|ON
81|46042.64583333334|-6.43671|forecast|46042|ON
82|46042.66666666666|-6|forecast|46042|OFF
83|46042.6875|-7.16779|forecast|46042|OFF
84|46042.70833333334|-7.85638|forecast|46042|OFF
85|46042.72916666666|-6|forecast|46042|OFF
86|46042.75|9.67714|forecast|46042|OFF
87|46042.77083333334|56.62538|forecast|46042|OFF
88|46042.79166666666|57.3|forecast|46042|OFF
89|46042.8125|73.20007|forecast|46042|OFF
90|46042.83333333334|73.20007|forecast|46042|OFF
91|46042.85416666666|64.8901|forecast|46042|OFF
92|46042.875|74.68738|forecast|46042|OFF
93|46042.89583333334|73.20007|forecast|46042|OFF
94|46042.91666666666|61.08817|forecast|46042|OFF
95|46042.9375|63.23909|forecast|46042|OFF
96|46042.95833333334|61.56018|forecast|46042|OFF
97|46042.97916666666|57.3|forecast|46042|OFF
"@

$newRowLines = $newRows -split "`n"
foreach ($line in $newRowLines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $p = $line.Split("|")
    $r = [int]$p[0]

    $wsDetailed.Cells.Item($r, 1).Value = [double]$p[1]
    $wsDetailed.Cells.Item($r, 1).NumberFormat = $dtFmt

    $wsDetailed.Cells.Item($r, 2).Value = [double]$p[2]

    $wsDetailed.Cells.Item($r, 3).Value = $p[3]

    $wsDetailed.Cells.Item($r, 4).Value = [double]$p[4]
    $wsDetailed.Cells.Item($r, 4).NumberFormat = $dFmt

    $wsDetailed.Cells.Item($r, 5).Value = $p[5]
}
